# Adds three new data sheets ("Eurocontrol", "WayPoint2050", "Destination2050")
# around the existing "Swiss" sheet, and populates header rows for the two
# new data sheets with the appropriate series labels.

$wb = $excel.ActiveWorkbook

# --- Eurocontrol sheet (placed before "Swiss") -----------------------------
$wsEuro = $wb.Worksheets.Add($wb.Worksheets.Item("Swiss"))
$wsEuro.Name = "Eurocontrol"

$euroHeaders = @(
    "Other (x)", "Other (y)",
    "SAF (x)", "SAF (y)",
    "ATM (x)", "ATM (y)",
    "Fleet revol (x)", "Fleet revol (y)",
    "Fleet evol (x)", "Fleet evol (y)"
)
for ($i = 0; $i -lt $euroHeaders.Length; $i++) {
    $wsEuro.Cells.Item(1, $i + 1).Value = $euroHeaders[$i]
}

# --- WayPoint2050 sheet (placed before "Swiss", after "Eurocontrol") -------
$wsWP = $wb.Worksheets.Add($wb.Worksheets.Item("Swiss"))
$wsWP.Name = "WayPoint2050"

$wpHeaders = @(
    "Market-Based Measure (x)", "Market-Based Measure (y)",
    "SAF (x)", "SAF (y)",
    "Operations and Infrastructure (x)", "Operations and Infrastructure (y)",
    "Technology (x)", "Technology (y)"
)
for ($i = 0; $i -lt $wpHeaders.Length; $i++) {
    $wsWP.Cells.Item(1, $i + 1).Value = $wpHeaders[$i]
}

# --- Destination2050 sheet (placed after "Swiss"), left empty --------------
$wsDest = $wb.Worksheets.Add($null, $wb.Worksheets.Item("Swiss"))
$wsDest.Name = "Destination2050"

# --- View / selection tweaks ------------------------------------------------

# Eurocontrol: zoom 125%, selection on A2
$wsEuro2 = $wb.Worksheets.Item("Eurocontrol")
$wsEuro2.Activate()
$excel.ActiveWindow.Zoom = 125
$wsEuro2.Range("A2").Select() | Out-Null

# WayPoint2050: default zoom, selection on G2
$wsWP2 = $wb.Worksheets.Item("WayPoint2050")
$wsWP2.Activate()
$wsWP2.Range("G2").Select() | Out-Null

# Swiss: selection moves to E25, and becomes the active tab
$wsSwiss2 = $wb.Worksheets.Item("Swiss")
$wsSwiss2.Activate()
$wsSwiss2.Range("E25").Select() | Out-Null
